# Weekly Fruta/Hortaliza price-sheet refresh (Vega Modelo de Temuco - Pera asiatica).
# Each existing record (rows 2-15) is rolled forward to the next week's figures,
# and a new row 16 is appended to keep the oldest record that rolled off row 15.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 45083
$ws.Cells.Item(2, 13).Value = 55
$ws.Cells.Item(2, 14).Value = 16000
$ws.Cells.Item(2, 15).Value = 16000
$ws.Cells.Item(2, 16).Value = 16000
$ws.Cells.Item(2, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(2, 19).Value = 889
$ws.Cells.Item(2, 20).Value = 18

# Row 3
$ws.Cells.Item(3, 4).Value = 44418
$ws.Cells.Item(3, 13).Value = 100
$ws.Cells.Item(3, 14).Value = 8000
$ws.Cells.Item(3, 15).Value = 8000
$ws.Cells.Item(3, 16).Value = 8000
$ws.Cells.Item(3, 19).Value = 533

# Row 4
$ws.Cells.Item(4, 4).Value = 44511
$ws.Cells.Item(4, 13).Value = 15
$ws.Cells.Item(4, 14).Value = 22000
$ws.Cells.Item(4, 15).Value = 22000
$ws.Cells.Item(4, 16).Value = 22000
$ws.Cells.Item(4, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(4, 19).Value = 1467
$ws.Cells.Item(4, 20).Value = 15

# Row 5
$ws.Cells.Item(5, 4).Value = 44217
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 55
$ws.Cells.Item(5, 14).Value = 18000
$ws.Cells.Item(5, 15).Value = 18000
$ws.Cells.Item(5, 16).Value = 18000
$ws.Cells.Item(5, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(5, 19).Value = 1000
$ws.Cells.Item(5, 20).Value = 18

# Row 6
$ws.Cells.Item(6, 4).Value = 44966
$ws.Cells.Item(6, 13).Value = 4
$ws.Cells.Item(6, 14).Value = 250000
$ws.Cells.Item(6, 15).Value = 250000
$ws.Cells.Item(6, 16).Value = 250000
$ws.Cells.Item(6, 17).Value = '$/bins (400 kilos)'
$ws.Cells.Item(6, 19).Value = 625
$ws.Cells.Item(6, 20).Value = 400

# Row 7
$ws.Cells.Item(7, 4).Value = 44966
$ws.Cells.Item(7, 13).Value = 80
$ws.Cells.Item(7, 14).Value = 15000
$ws.Cells.Item(7, 15).Value = 15000
$ws.Cells.Item(7, 16).Value = 15000
$ws.Cells.Item(7, 19).Value = 833

# Row 8
$ws.Cells.Item(8, 4).Value = 45079
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 14).Value = 18000
$ws.Cells.Item(8, 15).Value = 18000
$ws.Cells.Item(8, 16).Value = 18000
$ws.Cells.Item(8, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(8, 19).Value = 1000
$ws.Cells.Item(8, 20).Value = 18

# Row 9
$ws.Cells.Item(9, 4).Value = 44208
$ws.Cells.Item(9, 12).Value = 'Especial'
$ws.Cells.Item(9, 13).Value = 70
$ws.Cells.Item(9, 14).Value = 24000
$ws.Cells.Item(9, 15).Value = 24000
$ws.Cells.Item(9, 16).Value = 24000
$ws.Cells.Item(9, 19).Value = 1600

# Row 10
$ws.Cells.Item(10, 4).Value = 44495
$ws.Cells.Item(10, 13).Value = 50
$ws.Cells.Item(10, 14).Value = 24000
$ws.Cells.Item(10, 15).Value = 24000
$ws.Cells.Item(10, 16).Value = 24000
$ws.Cells.Item(10, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(10, 18).Value = 'China'
$ws.Cells.Item(10, 19).Value = 2400
$ws.Cells.Item(10, 20).Value = 10

# Row 11
$ws.Cells.Item(11, 4).Value = 44264
$ws.Cells.Item(11, 12).Value = 'Calibre 100'
$ws.Cells.Item(11, 13).Value = 50
$ws.Cells.Item(11, 14).Value = 20000
$ws.Cells.Item(11, 15).Value = 20000
$ws.Cells.Item(11, 16).Value = 20000
$ws.Cells.Item(11, 17).Value = '$/caja 18 kilos embalada'
$ws.Cells.Item(11, 19).Value = 1111

# Row 12
$ws.Cells.Item(12, 4).Value = 44427
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 55
$ws.Cells.Item(12, 14).Value = 7000
$ws.Cells.Item(12, 15).Value = 7000
$ws.Cells.Item(12, 16).Value = 7000
$ws.Cells.Item(12, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(12, 19).Value = 467
$ws.Cells.Item(12, 20).Value = 15

# Row 13
$ws.Cells.Item(13, 4).Value = 44601
$ws.Cells.Item(13, 13).Value = 30
$ws.Cells.Item(13, 14).Value = 28000
$ws.Cells.Item(13, 15).Value = 28000
$ws.Cells.Item(13, 16).Value = 28000
$ws.Cells.Item(13, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(13, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(13, 19).Value = 1556
$ws.Cells.Item(13, 20).Value = 18

# Row 14
$ws.Cells.Item(14, 4).Value = 44392
$ws.Cells.Item(14, 12).Value = 'Especial'
$ws.Cells.Item(14, 13).Value = 500
$ws.Cells.Item(14, 14).Value = 7000
$ws.Cells.Item(14, 15).Value = 7000
$ws.Cells.Item(14, 16).Value = 7000
$ws.Cells.Item(14, 17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(14, 19).Value = 875
$ws.Cells.Item(14, 20).Value = 8

# Row 15
$ws.Cells.Item(15, 4).Value = 45085
$ws.Cells.Item(15, 13).Value = 110
$ws.Cells.Item(15, 14).Value = 16000
$ws.Cells.Item(15, 15).Value = 16000
$ws.Cells.Item(15, 16).Value = 16000
$ws.Cells.Item(15, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(15, 19).Value = 889
$ws.Cells.Item(15, 20).Value = 18

# Row 16 (new record, carried over from former row 15)
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(16, 3).Value = 'La Araucanía'
$ws.Cells.Item(16, 4).Value = 44411
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = 'Fruta'
$ws.Cells.Item(16, 7).Value = 100104
$ws.Cells.Item(16, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(16, 9).Value = 100104005
$ws.Cells.Item(16, 10).Value = 'Pera asiática'
$ws.Cells.Item(16, 11).Value = 'Hosui'
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 210
$ws.Cells.Item(16, 14).Value = 8000
$ws.Cells.Item(16, 15).Value = 8000
$ws.Cells.Item(16, 16).Value = 8000
$ws.Cells.Item(16, 17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(16, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(16, 19).Value = 1000
$ws.Cells.Item(16, 20).Value = 8

# Match the date format used by the other rows in column D
$ws.Cells.Item(16, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

